$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 108 (shifts existing rows 108-170 down to 109-171,
# and grows the used range to A1:T171).
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new weekly record.
$ws.Range("A108").Value() = 11
$ws.Range("B108").Value() = "Vega Monumental Concepción"
$ws.Range("C108").Value() = "Bíobío"
$ws.Range("D108").Value() = Get-Date -Year 2023 -Month 6 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Range("E108").Value() = 8
$ws.Range("F108").Value() = "Fruta"
$ws.Range("G108").Value() = 100108
$ws.Range("H108").Value() = "Tropicales y subtropicales"
$ws.Range("I108").Value() = 100108002
$ws.Range("J108").Value() = "Mango"
$ws.Range("K108").Value() = "Sin especificar"
$ws.Range("L108").Value() = "Primera"
$ws.Range("M108").Value() = 200
$ws.Range("N108").Value() = 9500
$ws.Range("O108").Value() = 10000
$ws.Range("P108").Value() = 9750
$ws.Range("Q108").Value() = "$/bandeja 4 kilos"
$ws.Range("R108").Value() = "Perú"
$ws.Range("S108").Value() = 2438
$ws.Range("T108").Value() = 4
